$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "CasesTab" row's Neo4j query (column B, row 2) dropped the optional
# "Cohort" column from its RETURN clause.
$ws.Range("B2").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`nMATCH (samp:sample)-->(c) `n WHERE samp.specific_sample_pathology IN [`"Pulmonary Adenocarcinoma`"]  `nOPTIONAL MATCH (co:cohort)<-[*]-(c)`n  WITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

# sheetView: zoom changed from 55% to 40%, selection moved from B4 to D4,
# and the frozen/scrolled top-left cell reset back to A1.
$ws.Activate()
$excel.ActiveWindow.Zoom = 40
$ws.Range("D4").Select()
